# Update the stats after the latest release.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$lo = $ws.ListObjects.Item("Data")

# --- Correct a value in the last existing row (YAML count for 44993) ---
# The dependent calculated columns (Total, ∆Total) recompute automatically.
$ws.Range("K64").Value = 385

# --- Append a new day's row (44994) to the table ---
$lo.ListRows.Add() | Out-Null

# Copy formatting (number formats / styles) from the previous row down,
# same as Excel does when a table grows by typing into the row below it.
$ws.Range("A64:AJ64").Copy()
$ws.Range("A65:AJ65").PasteSpecial(-4122) | Out-Null

# Raw data values for the new row
$ws.Range("A65").Value = 44994
$ws.Range("B65").Value = 329
$ws.Range("C65").Value = 125
$ws.Range("D65").Value = 111
$ws.Range("E65").Value = 266
$ws.Range("F65").Value = 219
$ws.Range("G65").Value = 5511
$ws.Range("I65").Value = 6540
$ws.Range("J65").Value = 1922
$ws.Range("K65").Value = 385
$ws.Range("L65").Value = 285
$ws.Range("M65").Value = 114
$ws.Range("N65").Value = 60
$ws.Range("Q65").Value = 2031
$ws.Range("R65").Value = 4263
$ws.Range("S65").Value = 68367
$ws.Range("T65").Value = 47647
$ws.Range("U65").Value = 1
$ws.Range("V65").Value = 1
$ws.Range("W65").Value = 253
$ws.Range("Y65").Value = 0
$ws.Range("Z65").Value = 173
$ws.Range("AB65").Value = 156
$ws.Range("AC65").Value = 167
$ws.Range("AD65").Value = 5
$ws.Range("AE65").Value = 0
$ws.Range("AF65").Value = 377
$ws.Range("AG65").Value = 1093
$ws.Range("AH65").Value = 11
# AI65 (GH runs) intentionally left blank, matching the source row.

# Calculated (table formula) columns for the new row
$ws.Range("H65").Formula = "=Data[[#This Row],[LoC]]-G64"
$ws.Range("O65").Formula = "=SUM(Data[[#This Row],[Shell]:[Bash]])"
$ws.Range("P65").Formula = "=Data[[#This Row],[Total]]-O64"
$ws.Range("X65").Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"
$ws.Range("AA65").Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"
$ws.Range("AJ65").Formula = "=SUM(Data[[#This Row],[Running]:[GH runs]])"

# Match the author's final selection in the sheet
$ws.Range("AI65").Select() | Out-Null

Write-Output "Appended row 65 (44994) and corrected K64"
